$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.8415738951385362
$ws.Range("E2").Value = 0.8415738951385362
$ws.Range("D3").Value = 0.0002850719474319215
$ws.Range("E3").Value = 0.0002850719474319215
$ws.Range("D4").Value = 0.9388877107138294
$ws.Range("E4").Value = 0.9388877107138294
$ws.Range("D5").Value = 0.0009133948923708115
$ws.Range("E5").Value = 0.0009133948923708115
$ws.Range("D6").Value = 0.03007579956617298
$ws.Range("E6").Value = 0.03007579956617298
$ws.Range("D7").Value = 0.9977060769112752
$ws.Range("E7").Value = 0.002293923088724781
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.0001006203521423363
$ws.Range("E8").Value = 0.9998993796478577
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.0001474052275416579
$ws.Range("E9").Value = 0.9998525947724584
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.0000007220998930128472
$ws.Range("E10").Value = 0.999999277900107
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.000007615828992815823
$ws.Range("E11").Value = 0.9999923841710072
$ws.Range("F11").Value = 4.862440586090088
$ws.Range("G11").Value = 0.4
$ws.Range("D12").Value = 0.9581685755394228
$ws.Range("E12").Value = 0.9581685755394228
$ws.Range("D13").Value = 0.00003759235189215114
$ws.Range("E13").Value = 0.00003759235189215114
$ws.Range("D14").Value = 0.9878562164160773
$ws.Range("E14").Value = 0.9878562164160773
$ws.Range("D15").Value = 0.000069446419436554
$ws.Range("E15").Value = 0.000069446419436554
$ws.Range("D16").Value = 0.008555686559505121
$ws.Range("E16").Value = 0.008555686559505121
$ws.Range("D17").Value = 0.9999922195363712
$ws.Range("E17").Value = 0.000007780463628814793
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.000002344993683226224
$ws.Range("E18").Value = 0.9999976550063168
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.00002673337821742255
$ws.Range("E19").Value = 0.9999732666217825
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.0000004179592805526816
$ws.Range("E20").Value = 0.9999995820407195
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.0000007412130801613835
$ws.Range("E21").Value = 0.9999992587869199
$ws.Range("F21").Value = 5.988943576812744
$ws.Range("G21").Value = 0.4
